$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44312
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 330000
$ws.Range("O2").Value = 330000
$ws.Range("P2").Value = 330000
$ws.Range("Q2").Value = '$/bins (350 kilos)'
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 943
$ws.Range("T2").Value = 350

# Row 3
$ws.Range("D3").Value = 44195
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 200000
$ws.Range("O3").Value = 210000
$ws.Range("P3").Value = 206000
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 589

# Row 4
$ws.Range("D4").Value = 44193
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 150000
$ws.Range("O4").Value = 150000
$ws.Range("P4").Value = 150000
$ws.Range("S4").Value = 429

# Row 5
$ws.Range("D5").Value = 44363
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 200000
$ws.Range("O5").Value = 230000
$ws.Range("P5").Value = 215000
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 614

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("K6").Value = 'Start Ruby'
$ws.Range("M6").Value = 140
$ws.Range("N6").Value = 9800
$ws.Range("O6").Value = 9800
$ws.Range("P6").Value = 9800
$ws.Range("Q6").Value = '$/caja 14 kilos empedrada'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 700
$ws.Range("T6").Value = 14

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("M7").Value = 14
$ws.Range("N7").Value = 150000
$ws.Range("O7").Value = 160000
$ws.Range("P7").Value = 155000
$ws.Range("Q7").Value = '$/bins (350 kilos)'
$ws.Range("S7").Value = 443
$ws.Range("T7").Value = 350

# Row 8
$ws.Range("D8").Value = 44309
$ws.Range("M8").Value = 16
$ws.Range("N8").Value = 350000
$ws.Range("O8").Value = 350000
$ws.Range("P8").Value = 350000
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44298
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 15
$ws.Range("N9").Value = 450000
$ws.Range("O9").Value = 450000
$ws.Range("P9").Value = 450000
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 1286

# Row 10
$ws.Range("D10").Value = 44298
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 430000
$ws.Range("O10").Value = 430000
$ws.Range("P10").Value = 430000
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 1229

# Row 11
$ws.Range("D11").Value = 44189
$ws.Range("M11").Value = 16
$ws.Range("N11").Value = 150000
$ws.Range("O11").Value = 150000
$ws.Range("P11").Value = 150000
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 429

# Row 12
$ws.Range("D12").Value = 44196
$ws.Range("K12").Value = 'Red Blush'
$ws.Range("M12").Value = 12
$ws.Range("N12").Value = 130000
$ws.Range("O12").Value = 130000
$ws.Range("P12").Value = 130000
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 371

# Row 13
$ws.Range("D13").Value = 44308
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 280000
$ws.Range("O13").Value = 280000
$ws.Range("P13").Value = 280000
$ws.Range("S13").Value = 800

# Row 14
$ws.Range("D14").Value = 44208
$ws.Range("M14").Value = 16
$ws.Range("N14").Value = 180000
$ws.Range("O14").Value = 180000
$ws.Range("P14").Value = 180000
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 514

# Row 15
$ws.Range("D15").Value = 44376
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 180000
$ws.Range("O15").Value = 180000
$ws.Range("P15").Value = 180000
$ws.Range("R15").Value = 'Hijuelas'
$ws.Range("S15").Value = 514

# Row 16
$ws.Range("D16").Value = 44376
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 16
$ws.Range("N16").Value = 150000
$ws.Range("O16").Value = 150000
$ws.Range("P16").Value = 150000
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 429

# Row 17
$ws.Range("D17").Value = 44389
$ws.Range("M17").Value = 18
$ws.Range("R17").Value = 'Provincia de Quillota'

# Row 18
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 8
$ws.Range("N18").Value = 200000
$ws.Range("O18").Value = 200000
$ws.Range("P18").Value = 200000
$ws.Range("S18").Value = 571

# Row 19
$ws.Range("D19").Value = 44201
$ws.Range("M19").Value = 16
$ws.Range("N19").Value = 170000
$ws.Range("O19").Value = 170000
$ws.Range("P19").Value = 170000
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 486

# Row 20
$ws.Range("D20").Value = 44356
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 24
$ws.Range("O20").Value = 230000
$ws.Range("P20").Value = 215000
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 614

# Row 21
$ws.Range("D21").Value = 44400
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 140
$ws.Range("N21").Value = 9800
$ws.Range("O21").Value = 9800
$ws.Range("P21").Value = 9800
$ws.Range("Q21").Value = '$/caja 14 kilos empedrada'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 700
$ws.Range("T21").Value = 14
